# Apply scheduled-runner price/profit updates across all sheets
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H9").Value = 200
$ws.Range("I9").Value = 0
$ws.Range("J9").Value = 200
$ws.Range("K9").Value = 0
$ws.Range("L9").Value = 200
$ws.Range("N9").Value = -538
$ws.Range("H16").Value = 0
$ws.Range("I16").Value = 0
$ws.Range("J16").Value = 0
$ws.Range("K16").Value = 0
$ws.Range("L16").Value = 0
$ws.Range("H40").Value = 5610.722
$ws.Range("I40").Value = 3916.5833
$ws.Range("J40").Value = 8999
$ws.Range("K40").Value = 3916.5833
$ws.Range("L40").Value = 8999
$ws.Range("M40").Value = -3741.5833
$ws.Range("N40").Value = -9349
$ws.Range("H88").Value = 8794.833000000001
$ws.Range("I88").Value = 499
$ws.Range("J88").Value = 10454
$ws.Range("K88").Value = 499
$ws.Range("L88").Value = 10454
$ws.Range("M88").Value = -93
$ws.Range("N88").Value = -11266
$ws.Range("H91").Value = 8794.833000000001
$ws.Range("I91").Value = 499
$ws.Range("J91").Value = 10454
$ws.Range("K91").Value = 499
$ws.Range("L91").Value = 10454
$ws.Range("M91").Value = 905
$ws.Range("N91").Value = -13262
$ws.Range("H127").Value = 873.9091
$ws.Range("I127").Value = 624.8889
$ws.Range("J127").Value = 1994.5
$ws.Range("K127").Value = 1874.6667
$ws.Range("L127").Value = 5983.5
$ws.Range("M127").Value = 3085.3333
$ws.Range("N127").Value = -15903.5
$ws.Range("H132").Value = 46078.348
$ws.Range("I132").Value = 58290.723
$ws.Range("J132").Value = 2113.8
$ws.Range("K132").Value = 174872.169
$ws.Range("L132").Value = 6341.400000000001
$ws.Range("M132").Value = -172342.169
$ws.Range("N132").Value = -11401.4
$ws.Range("H138").Value = 7949.566
$ws.Range("J138").Value = 8566.333000000001
$ws.Range("L138").Value = 25698.999
$ws.Range("N138").Value = -35978.999
$ws.Range("M9").ClearContents()
$ws.Range("M16").ClearContents()
$ws.Range("N16").ClearContents()

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 2753.7646
$ws.Range("I32").Value = 2789.6736
$ws.Range("K32").Value = 2789.6736
$ws.Range("M32").Value = -2502.6736
$ws.Range("I74").Value = 4200.3125
$ws.Range("J74").Value = 2779989.2
$ws.Range("K74").Value = 4200.3125
$ws.Range("L74").Value = 2779989.2
$ws.Range("M74").Value = -3326.3125
$ws.Range("N74").Value = -2781737.2
$ws.Range("I77").Value = 4200.3125
$ws.Range("J77").Value = 2779989.2
$ws.Range("K77").Value = 21001.5625
$ws.Range("L77").Value = 13899946
$ws.Range("M77").Value = -16633.5625
$ws.Range("N77").Value = -13908682
$ws.Range("H102").Value = 615.7
$ws.Range("I102").Value = 678.4286
$ws.Range("J102").Value = 469.33334
$ws.Range("K102").Value = 678.4286
$ws.Range("L102").Value = 469.33334
$ws.Range("M102").Value = 943.5714
$ws.Range("N102").Value = -3713.33334
$ws.Range("H122").Value = 5258.9614
$ws.Range("I122").Value = 4769.909
$ws.Range("K122").Value = 14309.727
$ws.Range("M122").Value = -11859.727

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H23").Value = 11006.5
$ws.Range("J23").Value = 2000
$ws.Range("L23").Value = 2000
$ws.Range("N23").Value = -2566
$ws.Range("H43").Value = 0
$ws.Range("J43").Value = 0
$ws.Range("L43").Value = 0
$ws.Range("H86").Value = 2458.0908
$ws.Range("I86").Value = 1983.5
$ws.Range("J86").Value = 3027.6
$ws.Range("K86").Value = 1983.5
$ws.Range("L86").Value = 3027.6
$ws.Range("M86").Value = -860.5
$ws.Range("N86").Value = -5273.6
$ws.Range("H89").Value = 2458.0908
$ws.Range("I89").Value = 1983.5
$ws.Range("J89").Value = 3027.6
$ws.Range("K89").Value = 9917.5
$ws.Range("L89").Value = 15138
$ws.Range("M89").Value = -4301.5
$ws.Range("N89").Value = -26370
$ws.Range("H94").Value = 28802.25
$ws.Range("J94").Value = 2605
$ws.Range("L94").Value = 2605
$ws.Range("N94").Value = -3507
$ws.Range("H99").Value = 1689.091
$ws.Range("I99").Value = 1604.125
$ws.Range("J99").Value = 1915.6666
$ws.Range("K99").Value = 1604.125
$ws.Range("L99").Value = 1915.6666
$ws.Range("M99").Value = -106.125
$ws.Range("N99").Value = -4911.6666
$ws.Range("N43").ClearContents()

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 2743.0688
$ws.Range("J31").Value = 2426.1365
$ws.Range("L31").Value = 2426.1365
$ws.Range("N31").Value = -3016.1365
$ws.Range("H34").Value = 2743.0688
$ws.Range("J34").Value = 2426.1365
$ws.Range("L34").Value = 2426.1365
$ws.Range("N34").Value = -2830.1365
$ws.Range("H59").Value = 55095.832
$ws.Range("J59").Value = 60115
$ws.Range("L59").Value = 60115
$ws.Range("N59").Value = -62405
$ws.Range("H62").Value = 22596.6
$ws.Range("I62").Value = 19490
$ws.Range("K62").Value = 19490
$ws.Range("M62").Value = -18866
$ws.Range("H65").Value = 22596.6
$ws.Range("I65").Value = 19490
$ws.Range("K65").Value = 97450
$ws.Range("M65").Value = -94330
$ws.Range("H105").Value = 3985.2856
$ws.Range("I105").Value = 3985.2856
$ws.Range("K105").Value = 3985.2856
$ws.Range("M105").Value = -2238.2856
$ws.Range("H129").Value = 0
$ws.Range("I129").Value = 0
$ws.Range("J129").Value = 0
$ws.Range("K129").Value = 0
$ws.Range("L129").Value = 0
$ws.Range("H130").Value = 0
$ws.Range("I130").Value = 0
$ws.Range("J130").Value = 0
$ws.Range("K130").Value = 0
$ws.Range("L130").Value = 0
$ws.Range("H131").Value = 21728.8
$ws.Range("I131").Value = 0
$ws.Range("J131").Value = 21728.8
$ws.Range("K131").Value = 0
$ws.Range("L131").Value = 21728.8
$ws.Range("N131").Value = -31808.8
$ws.Range("H132").Value = 2133.7058
$ws.Range("I132").Value = 1805.5333
$ws.Range("J132").Value = 4595
$ws.Range("K132").Value = 5416.5999
$ws.Range("L132").Value = 13785
$ws.Range("M132").Value = -2886.5999
$ws.Range("N132").Value = -18845
$ws.Range("H133").Value = 0
$ws.Range("I133").Value = 0
$ws.Range("J133").Value = 0
$ws.Range("K133").Value = 0
$ws.Range("L133").Value = 0
$ws.Range("H134").Value = 0
$ws.Range("I134").Value = 0
$ws.Range("J134").Value = 0
$ws.Range("K134").Value = 0
$ws.Range("L134").Value = 0
$ws.Range("H135").Value = 55000
$ws.Range("I135").Value = 0
$ws.Range("J135").Value = 55000
$ws.Range("K135").Value = 0
$ws.Range("L135").Value = 55000
$ws.Range("N135").Value = -65140
$ws.Range("H137").Value = 0
$ws.Range("I137").Value = 0
$ws.Range("J137").Value = 0
$ws.Range("K137").Value = 0
$ws.Range("L137").Value = 0
$ws.Range("H138").Value = 0
$ws.Range("I138").Value = 0
$ws.Range("J138").Value = 0
$ws.Range("K138").Value = 0
$ws.Range("L138").Value = 0
$ws.Range("H139").Value = 0
$ws.Range("I139").Value = 0
$ws.Range("J139").Value = 0
$ws.Range("K139").Value = 0
$ws.Range("L139").Value = 0
$ws.Range("H140").Value = 0
$ws.Range("I140").Value = 0
$ws.Range("J140").Value = 0
$ws.Range("K140").Value = 0
$ws.Range("L140").Value = 0
$ws.Range("H141").Value = 55147.5
$ws.Range("I141").Value = 55147.5
$ws.Range("J141").Value = 0
$ws.Range("K141").Value = 55147.5
$ws.Range("L141").Value = 0
$ws.Range("M141").Value = -49967.5

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H2").Value = 1061.7778
$ws.Range("I2").Value = 1719.375
$ws.Range("J2").Value = 535.7
$ws.Range("K2").Value = 10316.25
$ws.Range("L2").Value = 3214.2
$ws.Range("M2").Value = -10203.25
$ws.Range("N2").Value = -3440.2
$ws.Range("H38").Value = 17.1
$ws.Range("I38").Value = 13.285714
$ws.Range("J38").Value = 26
$ws.Range("K38").Value = 39.857142
$ws.Range("L38").Value = 78
$ws.Range("M38").Value = 307.142858
$ws.Range("N38").Value = -772
$ws.Range("H54").Value = 60000
$ws.Range("I54").Value = 0
$ws.Range("J54").Value = 60000
$ws.Range("K54").Value = 0
$ws.Range("L54").Value = 180000
$ws.Range("N54").Value = -181118
$ws.Range("H113").Value = 353.15384
$ws.Range("I113").Value = 200.2
$ws.Range("J113").Value = 448.75
$ws.Range("K113").Value = 600.5999999999999
$ws.Range("L113").Value = 1346.25
$ws.Range("M113").Value = 1569.4
$ws.Range("N113").Value = -5686.25
$ws.Range("H137").Value = 377714
$ws.Range("J137").Value = 397486.84
$ws.Range("L137").Value = 1192460.52
$ws.Range("N137").Value = -1202660.52
$ws.Range("M54").ClearContents()

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 1448
$ws.Range("I80").Value = 0
$ws.Range("J80").Value = 1448
$ws.Range("K80").Value = 0
$ws.Range("L80").Value = 1448
$ws.Range("N80").Value = -3444
$ws.Range("H83").Value = 1448
$ws.Range("I83").Value = 0
$ws.Range("J83").Value = 1448
$ws.Range("K83").Value = 0
$ws.Range("L83").Value = 7240
$ws.Range("N83").Value = -17224
$ws.Range("H122").Value = 2845.1667
$ws.Range("I122").Value = 2870.611
$ws.Range("J122").Value = 2768.8333
$ws.Range("K122").Value = 8611.832999999999
$ws.Range("L122").Value = 8306.499899999999
$ws.Range("M122").Value = -6161.832999999999
$ws.Range("N122").Value = -13206.4999
$ws.Range("H132").Value = 8850.714
$ws.Range("I132").Value = 4395.4
$ws.Range("K132").Value = 13186.2
$ws.Range("M132").Value = -10656.2
$ws.Range("M80").ClearContents()
$ws.Range("M83").ClearContents()

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H102").Value = 18999.5
$ws.Range("J102").Value = 18999.5
$ws.Range("L102").Value = 18999.5
$ws.Range("N102").Value = -25489.5

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H15").Value = 38999.2
$ws.Range("H51").Value = 13356.667
$ws.Range("I51").Value = 10035
$ws.Range("J51").Value = 20000
$ws.Range("K51").Value = 10035
$ws.Range("L51").Value = 20000
$ws.Range("M51").Value = -9525
$ws.Range("N51").Value = -21020
$ws.Range("H107").Value = 1910.1666
$ws.Range("I107").Value = 1271
$ws.Range("K107").Value = 3813
$ws.Range("M107").Value = -1893
